$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new transaction rows right after the existing row 2 (header is row 1).
# This pushes the old rows 3:68 down to 5:70, matching the target layout.
$ws.Rows("3:4").Insert()

# New row 3: Deposit / Wiretransfer / 7576.5357999999997
$ws.Range("E3").Value = "Deposit"
$ws.Range("N3").Value = "Wiretransfer"
$ws.Range("T3").Value = 7576.5357999999997

# New row 4: Withdrawal / Credit Card / 269.6825
$ws.Range("E4").Value = "Withdrawal"
$ws.Range("N4").Value = "Credit Card"
$ws.Range("T4").Value = 269.6825

# Update the sheet selection/active cell as recorded after the edit.
$ws.Range("U1:AB1048576").Select()
